# Final restructuring before golive
#
# Appends three new "variant" rows (242-244) at the bottom of Sheet1,
# mirroring the "mensualidad / 12 dias / 8 dias" boarding-plan rows that
# already exist earlier in the sheet (rows 66-68), but filed under the
# newer product grouping (column C = 41) that rows 239-241 belong to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 242 - "mensualidad (de lunes a viernes)"
# ---------------------------------------------------------------------
$ws.Range("A242").Formula = "=A241+1"
$ws.Range("B242").Formula = "=B241"
$ws.Range("C242").Value   = 41
$ws.Range("D242").Value2  = $ws.Range("D66").Value2
$ws.Range("E242").Value   = 500000
$ws.Range("F242").Formula = "=E242*90%"
$ws.Range("G242").Formula = "=E242*5%"
$ws.Range("H242").Value   = 1
$ws.Range("I242").Value   = 1
$ws.Range("J242").Value   = 0
$ws.Range("K242").Value   = 30
$ws.Range("L242").Value2  = $ws.Range("L66").Value2
$ws.Range("M242").Value2  = $ws.Range("M66").Value2

# ---------------------------------------------------------------------
# Row 243 - "12 dias al mes (3 dias a la semana)"
# ---------------------------------------------------------------------
$ws.Range("A243").Formula = "=A242+1"
$ws.Range("B243").Formula = "=B242"
$ws.Range("C243").Value   = 41
$ws.Range("D243").Value2  = $ws.Range("D67").Value2
$ws.Range("E243").Value   = 300000
$ws.Range("F243").Formula = "=E243*90%"
$ws.Range("G243").Formula = "=E243*5%"
$ws.Range("H243").Value   = 1
$ws.Range("I243").Value   = 1
$ws.Range("J243").Value   = 0
$ws.Range("K243").Value   = 30
$ws.Range("L243").Formula = "=D243"
$ws.Range("M243").Value2  = $ws.Range("M67").Value2

# ---------------------------------------------------------------------
# Row 244 - "8 dias al mes (2 dias por semana)"
# ---------------------------------------------------------------------
$ws.Range("A244").Formula = "=A243+1"
$ws.Range("B244").Formula = "=B243"
$ws.Range("C244").Value   = 41
$ws.Range("D244").Value2  = $ws.Range("D68").Value2
$ws.Range("E244").Value   = 250000
$ws.Range("F244").Formula = "=E244*90%"
$ws.Range("G244").Formula = "=E244*5%"
$ws.Range("H244").Value   = 1
$ws.Range("I244").Value   = 1
$ws.Range("J244").Value   = 0
$ws.Range("K244").Value   = 30
$ws.Range("L244").Formula = "=D244"
$ws.Range("M244").Value2  = $ws.Range("M68").Value2

# ---------------------------------------------------------------------
# Carry over the cell formatting from the template rows (66-68) onto
# the newly written cells, column by column, without disturbing the
# values/formulas just written (PasteSpecial xlPasteFormats = -4122).
# ---------------------------------------------------------------------
$xlPasteFormats = -4122
$srcRows = @(66, 67, 68)
$dstRows = @(242, 243, 244)
$cols = @("A","B","D","E","F","G","H","I","J","K","L","M")

for ($i = 0; $i -lt 3; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]
    foreach ($col in $cols) {
        $ws.Range("$col$srcRow").Copy()
        $ws.Range("$col$dstRow").PasteSpecial($xlPasteFormats)
    }
}

# Column C: the product-grouping id now uses a dedicated Arial style
# (no explicit size/colour override) instead of the template's style.
foreach ($dstRow in $dstRows) {
    $ws.Range("C$dstRow").Font.Name = "Arial"
}
